# Apply the net-name/net-label reordering edits described by the diff,
# plus the row-21 height bump on the BoM sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "BoM"
# ---------------------------------------------------------------------
$bom = $wb.Worksheets.Item("BoM")

# Row 10 (C1 x47uF cap): Earth,Net-(U1-UCAP) -> Net-(U1-UCAP),Earth
$bom.Range("AC10").Value = "Net-(U1-UCAP),Earth"
$bom.Range("AD10").Value = "Net-(U1-UCAP),Earth"

# Row 12 (D1 diode): +5V,/RESET2 -> /RESET2,+5V ; RESET2 -> RESET2,+5V
$bom.Range("AC12").Value = "/RESET2,+5V"
$bom.Range("AD12").Value = "RESET2,+5V"

# Row 15 (J5 ICSP header)
$bom.Range("AC15").Value = "GND,/MISO2,/RESET2,/SCK2,/MOSI2,+5V"
$bom.Range("AD15").Value = "MOSI2,+5V"

# Row 16 (J3 connector)
$bom.Range("AC16").Value = "Net-(J3-Pin_3),Net-(J3-Pin_1),Net-(J3-Pin_2),Net-(J3-Pin_4),Net-(J3-Pin_5)"
$bom.Range("AD16").Value = "Net-(J3-Pin_3),Net-(J3-Pin_1),Net-(J3-Pin_2),Net-(J3-Pin_4),Net-(J3-Pin_5)"

# Row 17 (J6 connector)
$bom.Range("AC17").Value = "Net-(J6-Pin_1),Net-(J6-Pin_2),Net-(J6-Pin_3),Net-(J6-Pin_5),Net-(J6-Pin_6),Net-(J6-Pin_4)"
$bom.Range("AD17").Value = "Net-(J6-Pin_1),Net-(J6-Pin_2),Net-(J6-Pin_3),Net-(J6-Pin_5),Net-(J6-Pin_6),Net-(J6-Pin_4)"

# Row 20 (R4 resistor) shares the same underlying text as row 12
$bom.Range("AC20").Value = "/RESET2,+5V"
$bom.Range("AD20").Value = "RESET2,+5V"

# Row 21 (U1 MCU) - big net list, reordered
$bom.Range("AC21").Value = "GND,Net-(J3-Pin_3),Net-(J3-Pin_4),Net-(J6-Pin_5),Net-(J6-Pin_4),Net-(J4-Pin_3),VBUS,Net-(J4-Pin_1),Net-(J6-Pin_2),/TXLED,/RESET2,Net-(U1-D-),Net-(U1-PC0{slash}XTAL2),/MISO2,Net-(J6-Pin_3),unconnected-(U1-PB0-Pad14),/MOSI2,Net-(J4-Pin_4),Net-(J3-Pin_1),Net-(J3-Pin_2),+5V,Net-(J6-Pin_6),Net-(U1-D+),Earth,Net-(J4-Pin_2),/DTR,/SCK2,Net-(U1-XTAL1),Net-(J3-Pin_5),Net-(U1-UCAP),/RXLED"

# Row 21 height: 105 -> 120
$bom.Rows.Item(21).RowHeight = 120

# ---------------------------------------------------------------------
# Sheet "DNF"
# ---------------------------------------------------------------------
$dnf = $wb.Worksheets.Item("DNF")

# Row 9 (C3/C4 caps)
$dnf.Range("AC9").Value = "Net-(U1-XTAL1),GND"
$dnf.Range("AD9").Value = "Net-(U1-XTAL1),GND"

# Row 10 (F1 polyfuse)
$dnf.Range("AC10").Value = "VBUS,Net-(J2-VBUS)"
$dnf.Range("AD10").Value = "VBUS,Net-(J2-VBUS)"

# Row 12 (J4 connector)
$dnf.Range("AC12").Value = "Net-(J4-Pin_1),Net-(J4-Pin_3),Net-(J4-Pin_4),Net-(J4-Pin_2)"
$dnf.Range("AD12").Value = "Net-(J4-Pin_1),Net-(J4-Pin_3),Net-(J4-Pin_4),Net-(J4-Pin_2)"

# Row 13 (J2 USB B connector)
$dnf.Range("AC13").Value = "Net-(J2-D+),Net-(J2-VBUS),Net-(J2-D-),Net-(J2-Shield),Earth"
$dnf.Range("AD13").Value = "Net-(J2-D+),Net-(J2-VBUS),Net-(J2-D-),Net-(J2-Shield),Earth"
